$d = $word.ActiveDocument

# The document contains five "<id>...</id>" markers, each originally split
# across multiple runs (an opening "<id>" run, one or more runs holding the
# "p157r_N" identifier, and a closing "</id>" run). Collapse each one into a
# single run "<id>p157r_N</id>" that keeps the formatting of the opening
# "<id>" run (Courier New, color 7f6000, size 9pt).

for ($i = 1; $i -le 5; $i++) {
    $needle = "<id>p157r_" + $i + "</id>"
    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
